$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Table 1 - Caribbean warming")

# Header rename: HadISST 1987-2020 -> HadISST 1994-2020
$ws.Range("D1").Value = "HadISST..1994.2020."

# Data updates in column D (HadISST 1987/1994-2020)
$ws.Range("D2").Value = 0.2
$ws.Range("D3").Value = 0.54
$ws.Range("D4").Value = 0.17
$ws.Range("D5").Value = 0.46

# Data updates in column F (Pathfinder 1990-2019) for Caribbean Reefs rows
$ws.Range("F4").Value = 0.18
$ws.Range("F5").Value = 0.47
